$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the "Potencjał poprawy" formula across the whole column at once so
# Excel records it as a shared formula group (E2:E11), matching how the
# author re-entered/filled the formula down the table.
$ws.Range("E2:E11").Formula = "=D2*(100-C2)/100"

# Fill in the remaining "Column1" labels (Nr 5 / Nr 6 / Jak starczy czasu
# przed LEKiem) for the previously-blank rows. J9 is written before J8 so
# the new shared-string entries land in the same order as the target file
# (Nr 5 -> Nr 6 -> Jak starczy czasu przed LEKiem).
$ws.Range("J9").Value = "Nr 5"
$ws.Range("J8").Value = "Nr 6"
$ws.Range("J10").Value = "Jak starczy czasu przed LEKiem"

# Move/record the active selection as it was left after the edit.
$ws.Range("J18").Select() | Out-Null
